$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '92.221.18'
$ws.Range('E2').Value = '  -6.41%  '
$ws.Range('D3').Value = '3.318.38'
$ws.Range('E3').Value = '  -5.19%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.36'
$ws.Range('E5').Value = '  -10.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '617.12'
$ws.Range('E6').Value = '  -7.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.32'
$ws.Range('E7').Value = '  -11.09%  '
$ws.Range('E8').Value = '  -12.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.901'
$ws.Range('E10').Value = '  -14.81%  '
$ws.Range('D11').Value = '3.322.17'
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.190'
$ws.Range('E12').Value = '  -9.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.30'
$ws.Range('E13').Value = '  -13.93%  '
$ws.Range('D14').Value = '92.140.19'
$ws.Range('E14').Value = '  -6.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.82'
$ws.Range('E15').Value = '  -7.31%  '
$ws.Range('D16').Value = '3.935.71'
$ws.Range('E16').Value = '  -5.48%  '
$ws.Range('E17').Value = '  -8.71%  '
$ws.Range('E18').Value = '  -13.68%  '
$ws.Range('D19').Value = '3.315.11'
$ws.Range('E19').Value = '  -5.65%  '
$ws.Range('E20').Value = '  -13.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.70'
$ws.Range('E21').Value = '  -9.63%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '483.30'
$ws.Range('E22').Value = '  -8.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.435'
$ws.Range('E23').Value = '  -17.66%  '
$ws.Range('E24').Value = '  -10.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000180'
$ws.Range('E25').Value = '  -11.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.13'
$ws.Range('E26').Value = '  -10.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.41'
$ws.Range('E27').Value = '  -10.54%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.505.49'
$ws.Range('E28').Value = '  -4.95%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.22'
$ws.Range('E29').Value = '  -12.20%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.04'
$ws.Range('E31').Value = '  -11.78%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.56'
$ws.Range('E32').Value = '  -11.04%  '
$ws.Range('E33').Value = '  -12.03%  '
$ws.Range('B34').Value = 'Binance-PegBSC-USD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('B35').Value = 'Cronos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.167'
$ws.Range('E35').Value = '  -12.77%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '27.99'
$ws.Range('E36').Value = '  -11.04%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.514'
$ws.Range('E37').Value = '  -15.72%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '512.36'
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.25'
$ws.Range('E40').Value = '  -9.41%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.36'
$ws.Range('E41').Value = '  -11.72%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.145'
$ws.Range('E42').Value = '  -7.56%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.857'
$ws.Range('E43').Value = '  -8.13%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.00'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.65'
$ws.Range('E45').Value = '  -7.23%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.51'
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.31'
$ws.Range('E47').Value = '  -8.77%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.11'
$ws.Range('E48').Value = '  -6.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '51.91'
$ws.Range('E49').Value = '  -6.14%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0384'
$ws.Range('E50').Value = '  -12.30%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.04'
$ws.Range('E51').Value = '  -6.75%  '
